# Add a new employee row (Alex Jordan) to the Employee sheet, mirroring
# the existing rows: Firstname / Lastname / Username / Password(hyperlink).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "Alex"
$ws.Range("B6").Value = "Jordan"
$ws.Range("C6").Value = "alexJordan"
$ws.Range("D6").Value = "jordan@_2023!!!"

# Mirror the existing D5 "Password" cell: a mailto hyperlink styled like
# the workbook's other credential cells.
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:jordan@_2023!!!")
$ws.Range("D6").Style = "Hyperlink"

# Move the active selection, matching the saved session's last selection.
$null = $ws.Range("G7").Select()
